# Update "想去人数" (want-to-go count, column F) figures on the
# "展览" and "全部类型" sheets to the latest scraped snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 165
$ws1.Range("F5").Value  = 4732
$ws1.Range("F9").Value  = 480
$ws1.Range("F10").Value = 24
$ws1.Range("F11").Value = 17
$ws1.Range("F12").Value = 1343
$ws1.Range("F13").Value = 2937
$ws1.Range("F16").Value = 85
$ws1.Range("F18").Value = 2394
$ws1.Range("F21").Value = 31
$ws1.Range("F23").Value = 114
$ws1.Range("F24").Value = 48
$ws1.Range("F25").Value = 238

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 165
$ws4.Range("F6").Value  = 4732
$ws4.Range("F10").Value = 480
$ws4.Range("F11").Value = 24
$ws4.Range("F12").Value = 17
$ws4.Range("F13").Value = 1343
$ws4.Range("F14").Value = 2938
$ws4.Range("F17").Value = 85
$ws4.Range("F19").Value = 2394
$ws4.Range("F22").Value = 31
$ws4.Range("F24").Value = 114
$ws4.Range("F25").Value = 48
$ws4.Range("F26").Value = 238
